$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 238.88889
$ws.Range("I9").Value = 243.75
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 243.75
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = -74.75
$ws.Range("N9").Value = -538
$ws.Range("H18").Value = 975.5
$ws.Range("I18").Value = 966.6667
$ws.Range("K18").Value = 966.6667
$ws.Range("M18").Value = -682.6667
$ws.Range("H40").Value = 2541.125
$ws.Range("J40").Value = 2865.8
$ws.Range("L40").Value = 2865.8
$ws.Range("N40").Value = -3215.8
$ws.Range("H58").Value = 1972.2142
$ws.Range("I58").Value = 461.54544
$ws.Range("J58").Value = 7511.3335
$ws.Range("K58").Value = 1384.63632
$ws.Range("L58").Value = 22534.0005
$ws.Range("M58").Value = -1234.63632
$ws.Range("N58").Value = -22834.0005
$ws.Range("H69").Value = 324611
$ws.Range("J69").Value = 363624.88
$ws.Range("L69").Value = 1090874.64
$ws.Range("N69").Value = -1092622.64
$ws.Range("H70").Value = 7985.2856
$ws.Range("J70").Value = 5752.263
$ws.Range("L70").Value = 17256.789
$ws.Range("N70").Value = -17796.789
$ws.Range("H72").Value = 324611
$ws.Range("J72").Value = 363624.88
$ws.Range("L72").Value = 3272623.92
$ws.Range("N72").Value = -3281359.92
$ws.Range("H73").Value = 7985.2856
$ws.Range("J73").Value = 5752.263
$ws.Range("L73").Value = 17256.789
$ws.Range("N73").Value = -19128.789
$ws.Range("H74").Value = 3131.8
$ws.Range("I74").Value = 3131.8
$ws.Range("K74").Value = 3131.8
$ws.Range("M74").Value = -2195.8
$ws.Range("H76").Value = 4681
$ws.Range("I76").Value = 4183.2383
$ws.Range("K76").Value = 4183.2383
$ws.Range("M76").Value = -3868.2383
$ws.Range("H77").Value = 3131.8
$ws.Range("I77").Value = 3131.8
$ws.Range("K77").Value = 15659
$ws.Range("M77").Value = -10979
$ws.Range("H79").Value = 4681
$ws.Range("I79").Value = 4183.2383
$ws.Range("K79").Value = 4183.2383
$ws.Range("M79").Value = -3091.2383
$ws.Range("H86").Value = 2584.652
$ws.Range("I86").Value = 2439
$ws.Range("J86").Value = 2774
$ws.Range("K86").Value = 2439
$ws.Range("L86").Value = 2774
$ws.Range("M86").Value = -1316
$ws.Range("N86").Value = -5020
$ws.Range("H89").Value = 2584.652
$ws.Range("I89").Value = 2439
$ws.Range("J89").Value = 2774
$ws.Range("K89").Value = 12195
$ws.Range("L89").Value = 13870
$ws.Range("M89").Value = -6579
$ws.Range("N89").Value = -25102
$ws.Range("H94").Value = 7374.0835
$ws.Range("I94").Value = 8379.4
$ws.Range("J94").Value = 2347.5
$ws.Range("K94").Value = 8379.4
$ws.Range("L94").Value = 2347.5
$ws.Range("M94").Value = -7928.4
$ws.Range("N94").Value = -3249.5
$ws.Range("H112").Value = 7528.4473
$ws.Range("I112").Value = 1689
$ws.Range("J112").Value = 7686.2705
$ws.Range("K112").Value = 5067
$ws.Range("L112").Value = 23058.8115
$ws.Range("N112").Value = -25274.8115
$ws.Range("M112").Value = -3959
$ws.Range("H129").Value = 1199.4
$ws.Range("I129").Value = 799.6667
$ws.Range("K129").Value = 2399.0001
$ws.Range("M129").Value = 2600.9999
$ws.Range("H135").Value = 14027.195
$ws.Range("I135").Value = 1474.4166
$ws.Range("J135").Value = 104407.2
$ws.Range("K135").Value = 13269.7494
$ws.Range("L135").Value = 939664.7999999999
$ws.Range("M135").Value = -10734.7494
$ws.Range("N135").Value = -944734.7999999999
$ws.Range("H136").Value = 61000
$ws.Range("J136").Value = 61000
$ws.Range("L136").Value = 61000
$ws.Range("N136").Value = -71200
$ws.Range("H137").Value = 8140.7256
$ws.Range("I137").Value = 10110.257
$ws.Range("J137").Value = 3832.375
$ws.Range("K137").Value = 30330.771
$ws.Range("L137").Value = 11497.125
$ws.Range("M137").Value = -27780.771
$ws.Range("N137").Value = -16597.125
$ws.Range("H138").Value = 2868.0625
$ws.Range("I138").Value = 2160.3333
$ws.Range("J138").Value = 3778
$ws.Range("K138").Value = 6480.999899999999
$ws.Range("L138").Value = 11334
$ws.Range("M138").Value = -1340.999899999999
$ws.Range("N138").Value = -21614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1777.0588
$ws.Range("I2").Value = 1337.1364
$ws.Range("K2").Value = 1337.1364
$ws.Range("M2").Value = -1224.1364
$ws.Range("H32").Value = 17929.809
$ws.Range("I32").Value = 18209.521
$ws.Range("J32").Value = 8000
$ws.Range("K32").Value = 18209.521
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = -17922.521
$ws.Range("N32").Value = -8574
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H61").Value = 1037.5
$ws.Range("I61").Value = 1037.5
$ws.Range("K61").Value = 1037.5
$ws.Range("M61").Value = -825.5
$ws.Range("H74").Value = 240940.69
$ws.Range("I74").Value = 286567.66
$ws.Range("J74").Value = 1399
$ws.Range("K74").Value = 286567.66
$ws.Range("L74").Value = 1399
$ws.Range("M74").Value = -285693.66
$ws.Range("N74").Value = -3147
$ws.Range("H77").Value = 240940.69
$ws.Range("I77").Value = 286567.66
$ws.Range("J77").Value = 1399
$ws.Range("K77").Value = 1432838.3
$ws.Range("L77").Value = 6995
$ws.Range("M77").Value = -1428470.3
$ws.Range("N77").Value = -15731
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("H97").Value = 1929.8
$ws.Range("I97").Value = 1578.9286
$ws.Range("J97").Value = 2507.7058
$ws.Range("K97").Value = 1578.9286
$ws.Range("L97").Value = 2507.7058
$ws.Range("M97").Value = -1082.9286
$ws.Range("N97").Value = -3499.7058
$ws.Range("H116").Value = 1777.0588
$ws.Range("I116").Value = 1337.1364
$ws.Range("K116").Value = 1337.1364
$ws.Range("M116").Value = 956.8635999999999
$ws.Range("H132").Value = 2817.7407
$ws.Range("I132").Value = 1772.6666
$ws.Range("J132").Value = 3653.8
$ws.Range("K132").Value = 5317.9998
$ws.Range("L132").Value = 10961.4
$ws.Range("M132").Value = -2787.9998
$ws.Range("N132").Value = -16021.4
$ws.Range("H136").Value = 1037.5
$ws.Range("I136").Value = 1037.5
$ws.Range("K136").Value = 3112.5
$ws.Range("M136").Value = -562.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1777.0588
$ws.Range("I3").Value = 1337.1364
$ws.Range("K3").Value = 1337.1364
$ws.Range("M3").Value = -1223.1364
$ws.Range("H20").Value = 18473.896
$ws.Range("I20").Value = 25062.096
$ws.Range("J20").Value = 1179.875
$ws.Range("K20").Value = 25062.096
$ws.Range("L20").Value = 1179.875
$ws.Range("M20").Value = -24815.096
$ws.Range("N20").Value = -1673.875
$ws.Range("H105").Value = 4411.5713
$ws.Range("I105").Value = 2879.3684
$ws.Range("K105").Value = 2879.3684
$ws.Range("M105").Value = -1132.3684
$ws.Range("H107").Value = 52087.4
$ws.Range("I107").Value = 85141.664
$ws.Range("J107").Value = 2506
$ws.Range("K107").Value = 85141.664
$ws.Range("L107").Value = 2506
$ws.Range("M107").Value = -83221.664
$ws.Range("N107").Value = -6346
$ws.Range("H125").Value = 80888
$ws.Range("J125").Value = 80888
$ws.Range("L125").Value = 80888
$ws.Range("N125").Value = -90728
$ws.Range("H134").Value = 2837.9048
$ws.Range("I134").Value = 2720.4358
$ws.Range("K134").Value = 8161.307400000001
$ws.Range("M134").Value = -5626.307400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3707227.8
$ws.Range("I31").Value = 4350685
$ws.Range("J31").Value = 7349.75
$ws.Range("K31").Value = 4350685
$ws.Range("L31").Value = 7349.75
$ws.Range("M31").Value = -4350390
$ws.Range("N31").Value = -7939.75
$ws.Range("H34").Value = 3707227.8
$ws.Range("I34").Value = 4350685
$ws.Range("J34").Value = 7349.75
$ws.Range("K34").Value = 4350685
$ws.Range("L34").Value = 7349.75
$ws.Range("M34").Value = -4350483
$ws.Range("N34").Value = -7753.75
$ws.Range("H51").Value = 30998.309
$ws.Range("J51").Value = 30998.309
$ws.Range("L51").Value = 30998.309
$ws.Range("N51").Value = -32470.309
$ws.Range("H55").Value = 54249.5
$ws.Range("I55").Value = 29999
$ws.Range("K55").Value = 29999
$ws.Range("M55").Value = -29684
$ws.Range("H61").Value = 30998.309
$ws.Range("J61").Value = 30998.309
$ws.Range("L61").Value = 30998.309
$ws.Range("N61").Value = -31694.309
$ws.Range("H105").Value = 3912.5
$ws.Range("I105").Value = 3912.5
$ws.Range("K105").Value = 3912.5
$ws.Range("M105").Value = -2165.5
$ws.Range("H106").Value = 36000
$ws.Range("J106").Value = 36000
$ws.Range("L106").Value = 36000
$ws.Range("N106").Value = -38524
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("H132").Value = 36971.766
$ws.Range("I132").Value = 38970.062
$ws.Range("K132").Value = 116910.186
$ws.Range("M132").Value = -114380.186
$ws.Range("H134").Value = 2510.3896
$ws.Range("I134").Value = 1770.5483
$ws.Range("K134").Value = 5311.644899999999
$ws.Range("M134").Value = -2776.644899999999
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 49955628
$ws.Range("I4").Value = 32397842
$ws.Range("K4").Value = 97193526
$ws.Range("M4").Value = -97193414
$ws.Range("H11").Value = 3606175.8
$ws.Range("I11").Value = 4260707.5
$ws.Range("K11").Value = 12782122.5
$ws.Range("M11").Value = -12781982.5
$ws.Range("H12").Value = 199
$ws.Range("J12").Value = 178.38461
$ws.Range("L12").Value = 535.15383
$ws.Range("N12").Value = -881.15383
$ws.Range("H13").Value = 79.8
$ws.Range("I13").Value = 171
$ws.Range("J13").Value = 19
$ws.Range("K13").Value = 513
$ws.Range("L13").Value = 57
$ws.Range("M13").Value = -345
$ws.Range("N13").Value = -393
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("H37").Value = 39990
$ws.Range("J37").Value = 39990
$ws.Range("L37").Value = 119970
$ws.Range("N37").Value = -120194
$ws.Range("H97").Value = 2101.5
$ws.Range("J97").Value = 3970
$ws.Range("L97").Value = 11910
$ws.Range("N97").Value = -12902
$ws.Range("H113").Value = 1218
$ws.Range("I113").Value = 667
$ws.Range("K113").Value = 2001
$ws.Range("M113").Value = 169
$ws.Range("H131").Value = 186286.1
$ws.Range("J131").Value = 2014.5714
$ws.Range("L131").Value = 6043.7142
$ws.Range("N131").Value = -16123.7142
$ws.Range("H140").Value = 5600
$ws.Range("I140").Value = 5600
$ws.Range("K140").Value = 16800
$ws.Range("M140").Value = -11620
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 54
$ws.Range("I2").Value = 57.57143
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 57.57143
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 55.42857
$ws.Range("N2").Value = -230
$ws.Range("H29").Value = 887.4286
$ws.Range("I29").Value = 875
$ws.Range("J29").Value = 933
$ws.Range("K29").Value = 875
$ws.Range("L29").Value = 933
$ws.Range("M29").Value = -585
$ws.Range("N29").Value = -1513
$ws.Range("H31").Value = 2486.6
$ws.Range("I31").Value = 1366.6666
$ws.Range("J31").Value = 4166.5
$ws.Range("K31").Value = 1366.6666
$ws.Range("L31").Value = 4166.5
$ws.Range("M31").Value = -1074.6666
$ws.Range("N31").Value = -4750.5
$ws.Range("H37").Value = 2486.6
$ws.Range("I37").Value = 1366.6666
$ws.Range("J37").Value = 4166.5
$ws.Range("K37").Value = 1366.6666
$ws.Range("L37").Value = 4166.5
$ws.Range("M37").Value = -1089.6666
$ws.Range("N37").Value = -4720.5
$ws.Range("H70").Value = 6108.222
$ws.Range("I70").Value = 5971
$ws.Range("K70").Value = 5971
$ws.Range("M70").Value = -5701
$ws.Range("H73").Value = 6108.222
$ws.Range("I73").Value = 5971
$ws.Range("K73").Value = 5971
$ws.Range("M73").Value = -5035
$ws.Range("H80").Value = 6759.684
$ws.Range("I80").Value = 3916.2856
$ws.Range("J80").Value = 14721.2
$ws.Range("K80").Value = 3916.2856
$ws.Range("L80").Value = 14721.2
$ws.Range("M80").Value = -2918.2856
$ws.Range("N80").Value = -16717.2
$ws.Range("H83").Value = 6759.684
$ws.Range("I83").Value = 3916.2856
$ws.Range("J83").Value = 14721.2
$ws.Range("K83").Value = 19581.428
$ws.Range("L83").Value = 73606
$ws.Range("M83").Value = -14589.428
$ws.Range("N83").Value = -83590
$ws.Range("H113").Value = 2328.4
$ws.Range("I113").Value = 2160.75
$ws.Range("J113").Value = 2999
$ws.Range("K113").Value = 2160.75
$ws.Range("L113").Value = 2999
$ws.Range("M113").Value = 9.25
$ws.Range("N113").Value = -7339
$ws.Range("H117").Value = 35000
$ws.Range("J117").Value = 35000
$ws.Range("L117").Value = 35000
$ws.Range("N117").Value = -41884
$ws.Range("H126").Value = 2726.6155
$ws.Range("J126").Value = 3050.6667
$ws.Range("L126").Value = 9152.000100000001
$ws.Range("N126").Value = -14092.0001
$ws.Range("H136").Value = 12203.765
$ws.Range("J136").Value = 12203.765
$ws.Range("L136").Value = 36611.295
$ws.Range("N136").Value = -41711.295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3109.2727
$ws.Range("I40").Value = 3189.7
$ws.Range("K40").Value = 3189.7
$ws.Range("M40").Value = -3053.7
$ws.Range("H46").Value = 6273.125
$ws.Range("I46").Value = 2195.6667
$ws.Range("J46").Value = 7214.077
$ws.Range("K46").Value = 2195.6667
$ws.Range("L46").Value = 7214.077
$ws.Range("N46").Value = -7590.077
$ws.Range("M46").Value = -2007.6667
$ws.Range("H61").Value = 2301.1667
$ws.Range("I61").Value = 2261.4
$ws.Range("K61").Value = 2261.4
$ws.Range("M61").Value = -2059.4
$ws.Range("H113").Value = 2301.1667
$ws.Range("I113").Value = 2261.4
$ws.Range("K113").Value = 2261.4
$ws.Range("M113").Value = -91.40000000000009
$ws.Range("H132").Value = 4498.4116
$ws.Range("I132").Value = 3634.4546
$ws.Range("K132").Value = 10903.3638
$ws.Range("M132").Value = -8373.363799999999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H136").Value = 7004
$ws.Range("I136").Value = 7004
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 21012
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -18462
$ws.Range("N133").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("H126").Value = 717184
$ws.Range("I126").Value = 3522
$ws.Range("K126").Value = 10566
$ws.Range("M126").Value = -8096
$ws.Range("H132").Value = 4273.357
$ws.Range("J132").Value = 4298.5
$ws.Range("L132").Value = 12895.5
$ws.Range("N132").Value = -17955.5
$ws.Range("H136").Value = 38883.74
$ws.Range("I136").Value = 47152.8
$ws.Range("J136").Value = 7874.75
$ws.Range("K136").Value = 141458.4
$ws.Range("L136").Value = 23624.25
$ws.Range("M136").Value = -138908.4
$ws.Range("N136").Value = -28724.25
$ws.Range("N112").ClearContents()
